$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 13333
$ws.Range("I54").Value = 9999.5
$ws.Range("K54").Value = 9999.5
$ws.Range("M54").Value = -9513.5
$ws.Range("H132").Value = 228405.23
$ws.Range("I132").Value = 246883.7
$ws.Range("J132").Value = 3584
$ws.Range("K132").Value = 740651.1000000001
$ws.Range("L132").Value = 10752
$ws.Range("M132").Value = -738121.1000000001
$ws.Range("N132").Value = -15812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4282.57
$ws.Range("I32").Value = 2944.5227
$ws.Range("J32").Value = 14094.917
$ws.Range("K32").Value = 2944.5227
$ws.Range("L32").Value = 14094.917
$ws.Range("M32").Value = -2657.5227
$ws.Range("N32").Value = -14668.917
$ws.Range("H62").Value = 100249
$ws.Range("J62").Value = 100249
$ws.Range("L62").Value = 100249
$ws.Range("N62").Value = -101497
$ws.Range("H65").Value = 100249
$ws.Range("J65").Value = 100249
$ws.Range("L65").Value = 300747
$ws.Range("N65").Value = -306987
$ws.Range("H74").Value = 3666.0938
$ws.Range("I74").Value = 3102.4902
$ws.Range("K74").Value = 3102.4902
$ws.Range("M74").Value = -2228.4902
$ws.Range("H77").Value = 3666.0938
$ws.Range("I77").Value = 3102.4902
$ws.Range("K77").Value = 15512.451
$ws.Range("M77").Value = -11144.451
$ws.Range("H110").Value = 3204.4
$ws.Range("I110").Value = 2551.3076
$ws.Range("K110").Value = 2551.3076
$ws.Range("M110").Value = -506.3076000000001
$ws.Range("H111").Value = 75429.336
$ws.Range("J111").Value = 75429.336
$ws.Range("L111").Value = 75429.336
$ws.Range("N111").Value = -83609.336
$ws.Range("H122").Value = 2248.6365
$ws.Range("I122").Value = 2248.6365
$ws.Range("K122").Value = 6745.9095
$ws.Range("M122").Value = -4295.9095
$ws.Range("H132").Value = 4216.4336
$ws.Range("I132").Value = 2663.5283
$ws.Range("K132").Value = 7990.5849
$ws.Range("M132").Value = -5460.5849
$ws.Range("H134").Value = 83533.164
$ws.Range("J134").Value = 83533.164
$ws.Range("L134").Value = 83533.164
$ws.Range("N134").Value = -93673.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 172775
$ws.Range("J130").Value = 172775
$ws.Range("L130").Value = 172775
$ws.Range("N130").Value = -182815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = ""
$ws.Range("H122").Value = 1729.7273
$ws.Range("I122").Value = 1701.1666
$ws.Range("J122").Value = 1764
$ws.Range("K122").Value = 5103.4998
$ws.Range("L122").Value = 5292
$ws.Range("M122").Value = -2653.4998
$ws.Range("N122").Value = -10192
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 3379.4075
$ws.Range("I132").Value = 1754.5652
$ws.Range("K132").Value = 5263.6956
$ws.Range("M132").Value = -2733.6956
$ws.Range("H134").Value = 9818.878000000001
$ws.Range("I134").Value = 9818.878000000001
$ws.Range("K134").Value = 29456.634
$ws.Range("M134").Value = -26921.634
$ws.Range("H135").Value = 98995
$ws.Range("J135").Value = 98995
$ws.Range("L135").Value = 98995
$ws.Range("N135").Value = -109135
$ws.Range("H141").Value = 258702.8
$ws.Range("J141").Value = 291557.38
$ws.Range("L141").Value = 291557.38
$ws.Range("N141").Value = -301917.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 101249.75
$ws.Range("J37").Value = 101249.75
$ws.Range("L37").Value = 303749.25
$ws.Range("N37").Value = -303973.25
$ws.Range("H57").Value = 9502.091
$ws.Range("J57").Value = 9713.777
$ws.Range("L57").Value = 29141.331
$ws.Range("N57").Value = -30259.331
$ws.Range("H75").Value = 3447.9443
$ws.Range("I75").Value = 1508.25
$ws.Range("J75").Value = 4002.1428
$ws.Range("K75").Value = 4524.75
$ws.Range("L75").Value = 12006.4284
$ws.Range("M75").Value = -3526.75
$ws.Range("N75").Value = -14002.4284
$ws.Range("H78").Value = 3447.9443
$ws.Range("I78").Value = 1508.25
$ws.Range("J78").Value = 4002.1428
$ws.Range("K78").Value = 13574.25
$ws.Range("L78").Value = 36019.2852
$ws.Range("M78").Value = -8582.25
$ws.Range("N78").Value = -46003.2852
$ws.Range("H81").Value = 83338210
$ws.Range("I81").Value = 333333630
$ws.Range("J81").Value = 6394.8887
$ws.Range("K81").Value = 1000000890
$ws.Range("L81").Value = 19184.6661
$ws.Range("M81").Value = -999999767
$ws.Range("N81").Value = -21430.6661
$ws.Range("H82").Value = 8654.666999999999
$ws.Range("H84").Value = 83338210
$ws.Range("I84").Value = 333333630
$ws.Range("J84").Value = 6394.8887
$ws.Range("K84").Value = 3000002670
$ws.Range("L84").Value = 57553.99830000001
$ws.Range("M84").Value = -2999997054
$ws.Range("N84").Value = -68785.99830000001
$ws.Range("H85").Value = 8654.666999999999
$ws.Range("H86").Value = 1091.3334
$ws.Range("J86").Value = 1124.4
$ws.Range("L86").Value = 3373.2
$ws.Range("N86").Value = -5745.200000000001
$ws.Range("H89").Value = 1091.3334
$ws.Range("J89").Value = 1124.4
$ws.Range("L89").Value = 10119.6
$ws.Range("N89").Value = -21975.6
$ws.Range("H92").Value = 345.69232
$ws.Range("I92").Value = 271.5
$ws.Range("J92").Value = 464.4
$ws.Range("K92").Value = 814.5
$ws.Range("L92").Value = 1393.2
$ws.Range("M92").Value = 433.5
$ws.Range("N92").Value = -3889.2
$ws.Range("H131").Value = 8533.303
$ws.Range("J131").Value = 10255
$ws.Range("L131").Value = 30765
$ws.Range("N131").Value = -40845
$ws.Range("H133").Value = 3031.5454
$ws.Range("I133").Value = 3031.5454
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 9094.636200000001
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
$ws.Range("M133").Value = -4034.636200000001
$ws.Range("H134").Value = 931.25
$ws.Range("I134").Value = 931.25
$ws.Range("K134").Value = 2793.75
$ws.Range("M134").Value = 2276.25
$ws.Range("H137").Value = 8117.4
$ws.Range("I137").Value = 5997
$ws.Range("J137").Value = 9531
$ws.Range("K137").Value = 17991
$ws.Range("L137").Value = 28593
$ws.Range("M137").Value = -12891
$ws.Range("N137").Value = -38793
$ws.Range("H139").Value = 908.5714
$ws.Range("I139").Value = 908.5714
$ws.Range("K139").Value = 2725.7142
$ws.Range("M139").Value = 2414.2858
$ws.Range("H141").Value = 1232.375
$ws.Range("I141").Value = 1232.375
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3697.125
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1482.875
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2997.3333
$ws.Range("I126").Value = 2603.077
$ws.Range("J126").Value = 3638
$ws.Range("K126").Value = 7809.231000000001
$ws.Range("L126").Value = 10914
$ws.Range("M126").Value = -5339.231000000001
$ws.Range("N126").Value = -15854
$ws.Range("H132").Value = 2557.5806
$ws.Range("I132").Value = 2073.875
$ws.Range("J132").Value = 4216
$ws.Range("K132").Value = 6221.625
$ws.Range("L132").Value = 12648
$ws.Range("M132").Value = -3691.625
$ws.Range("N132").Value = -17708

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3863.4
$ws.Range("I7").Value = 3828
$ws.Range("K7").Value = 3828
$ws.Range("M7").Value = -3716
$ws.Range("H22").Value = 55531.95
$ws.Range("I22").Value = 144745.42
$ws.Range("K22").Value = 144745.42
$ws.Range("M22").Value = -144450.42
$ws.Range("H27").Value = 55531.95
$ws.Range("I27").Value = 144745.42
$ws.Range("K27").Value = 144745.42
$ws.Range("M27").Value = -144638.42
$ws.Range("H126").Value = 3863.4
$ws.Range("I126").Value = 3828
$ws.Range("K126").Value = 11484
$ws.Range("M126").Value = -9014
$ws.Range("H132").Value = 3994.6736
$ws.Range("I132").Value = 3289.8647
$ws.Range("K132").Value = 9869.5941
$ws.Range("M132").Value = -7339.5941
$ws.Range("H133").Value = 149500
$ws.Range("J133").Value = 149500
$ws.Range("L133").Value = 149500
$ws.Range("N133").Value = -154560
$ws.Range("H139").Value = 97853.75
$ws.Range("J139").Value = 97853.75
$ws.Range("L139").Value = 97853.75
$ws.Range("N139").Value = -108133.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74600.664
$ws.Range("J46").Value = 74925.75
$ws.Range("L46").Value = 74925.75
$ws.Range("N46").Value = -75387.75
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31262
$ws.Range("H107").Value = 4180.1665
$ws.Range("I107").Value = 1118.3684
$ws.Range("J107").Value = 9468.727999999999
$ws.Range("K107").Value = 3355.1052
$ws.Range("L107").Value = 28406.184
$ws.Range("M107").Value = -1435.1052
$ws.Range("N107").Value = -32246.184
$ws.Range("H132").Value = 11375291
$ws.Range("I132").Value = 1255.6061
$ws.Range("J132").Value = 38185520
$ws.Range("K132").Value = 3766.8183
$ws.Range("L132").Value = 114556560
$ws.Range("M132").Value = -1236.8183
$ws.Range("N132").Value = -114561620
$ws.Range("H134").Value = 74600.664
$ws.Range("J134").Value = 74925.75
$ws.Range("L134").Value = 224777.25
$ws.Range("N134").Value = -229847.25
